# Trade #10 closed at 2026-02-16 21:54:01 - leadlag UP +0.000%
#
# Appends the new trade row to both the "All Trades" sheet (next free row)
# and the strategy-specific "leadlag" sheet (next free row).

$wb = $excel.ActiveWorkbook

$tradeNum    = 10
$tradeDate   = "2026-02-16"
$tradeTime   = "21:54:01"
$strategy    = "leadlag"
$side        = "UP"
$entryPrice  = 68336.28999999999
$status      = "OPEN"
$pnlPct      = 0
$pnlDollar   = 0
$capitalAfter = 100
$confidence  = 0.75
$entryReason = "Binance leading with 0.089% move"

function Add-TradeRow {
    param($ws, $row)

    # Force the date/time columns to stay plain text instead of being
    # auto-converted into Excel date/time serials.
    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("C$row").NumberFormat = "@"

    $ws.Cells.Item($row, 1).Value  = $tradeNum     # A: Trade #
    $ws.Range("B$row").Value       = $tradeDate    # B: Date
    $ws.Range("C$row").Value       = $tradeTime    # C: Time
    $ws.Cells.Item($row, 4).Value  = $strategy      # D: Strategy
    $ws.Cells.Item($row, 5).Value  = $side           # E: Side
    $ws.Cells.Item($row, 6).Value  = $entryPrice     # F: Entry Price
    # G: Exit Price -- left blank (trade is still OPEN)
    $ws.Cells.Item($row, 8).Value  = $status         # H: Status
    $ws.Cells.Item($row, 9).Value  = $pnlPct         # I: P&L %
    $ws.Cells.Item($row, 10).Value = $pnlDollar     # J: P&L $
    $ws.Cells.Item($row, 11).Value = $capitalAfter  # K: Capital After
    $ws.Cells.Item($row, 12).Value = $confidence    # L: Confidence
    $ws.Cells.Item($row, 13).Value = $entryReason   # M: Entry Reason
    # N: Exit Reason -- left blank (trade is still OPEN)
    $ws.Cells.Item($row, 15).Value = 0              # O: Duration (min)
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades 11

$leadlag = $wb.Worksheets.Item("leadlag")
Add-TradeRow $leadlag 10
